$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 904.7826
$ws.Range("J129").Value = 1011.62067
$ws.Range("L129").Value = 3034.86201
$ws.Range("N129").Value = -13034.86201
$ws.Range("H137").Value = 1924415.8
$ws.Range("I137").Value = 4546469.5
$ws.Range("J137").Value = 1576.5
$ws.Range("K137").Value = 13639408.5
$ws.Range("L137").Value = 4729.5
$ws.Range("M137").Value = -13636858.5
$ws.Range("N137").Value = -9829.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 41751024
$ws.Range("I61").Value = 50050950
$ws.Range("J61").Value = 251388.75
$ws.Range("K61").Value = 50050950
$ws.Range("L61").Value = 251388.75
$ws.Range("M61").Value = -50050738
$ws.Range("N61").Value = -251812.75
$ws.Range("H74").Value = 15751360
$ws.Range("I74").Value = 17930026
$ws.Range("J74").Value = 500700
$ws.Range("K74").Value = 17930026
$ws.Range("L74").Value = 500700
$ws.Range("M74").Value = -17929152
$ws.Range("N74").Value = -502448
$ws.Range("H77").Value = 15751360
$ws.Range("I77").Value = 17930026
$ws.Range("J77").Value = 500700
$ws.Range("K77").Value = 89650130
$ws.Range("L77").Value = 2503500
$ws.Range("M77").Value = -89645762
$ws.Range("N77").Value = -2512236
$ws.Range("H132").Value = 49852.74
$ws.Range("I132").Value = 39078.816
$ws.Range("J132").Value = 69245.8
$ws.Range("K132").Value = 117236.448
$ws.Range("L132").Value = 207737.4
$ws.Range("M132").Value = -114706.448
$ws.Range("N132").Value = -212797.4
$ws.Range("H136").Value = 41751024
$ws.Range("I136").Value = 50050950
$ws.Range("J136").Value = 251388.75
$ws.Range("K136").Value = 150152850
$ws.Range("L136").Value = 754166.25
$ws.Range("M136").Value = -150150300
$ws.Range("N136").Value = -759266.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 14664
$ws.Range("I86").Value = 29838.25
$ws.Range("J86").Value = 2524.6
$ws.Range("K86").Value = 29838.25
$ws.Range("L86").Value = 2524.6
$ws.Range("M86").Value = -28715.25
$ws.Range("N86").Value = -4770.6
$ws.Range("H89").Value = 14664
$ws.Range("I89").Value = 29838.25
$ws.Range("J89").Value = 2524.6
$ws.Range("K89").Value = 149191.25
$ws.Range("L89").Value = 12623
$ws.Range("M89").Value = -143575.25
$ws.Range("N89").Value = -23855
$ws.Range("H134").Value = 1588.76
$ws.Range("I134").Value = 1705.75
$ws.Range("K134").Value = 5117.25
$ws.Range("M134").Value = -2582.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2511.8096
$ws.Range("I31").Value = 1381.3125
$ws.Range("J31").Value = 6129.4
$ws.Range("K31").Value = 1381.3125
$ws.Range("L31").Value = 6129.4
$ws.Range("M31").Value = -1086.3125
$ws.Range("N31").Value = -6719.4
$ws.Range("H34").Value = 2511.8096
$ws.Range("I34").Value = 1381.3125
$ws.Range("J34").Value = 6129.4
$ws.Range("K34").Value = 1381.3125
$ws.Range("L34").Value = 6129.4
$ws.Range("M34").Value = -1179.3125
$ws.Range("N34").Value = -6533.4
$ws.Range("H58").Value = 23257478
$ws.Range("I58").Value = 35715410
$ws.Range("J58").Value = 2674.6
$ws.Range("K58").Value = 35715410
$ws.Range("L58").Value = 2674.6
$ws.Range("M58").Value = -35715207
$ws.Range("N58").Value = -3080.6
$ws.Range("H132").Value = 40895.81
$ws.Range("I132").Value = 2506.8
$ws.Range("J132").Value = 93244.45
$ws.Range("K132").Value = 7520.400000000001
$ws.Range("L132").Value = 279733.35
$ws.Range("M132").Value = -4990.400000000001
$ws.Range("N132").Value = -284793.35
$ws.Range("H134").Value = 41948.668
$ws.Range("I134").Value = 2645.4666
$ws.Range("J134").Value = 91077.664
$ws.Range("K134").Value = 7936.399800000001
$ws.Range("L134").Value = 273232.992
$ws.Range("M134").Value = -5401.399800000001
$ws.Range("N134").Value = -278302.992
$ws.Range("H136").Value = 23257478
$ws.Range("I136").Value = 35715410
$ws.Range("J136").Value = 2674.6
$ws.Range("K136").Value = 107146230
$ws.Range("L136").Value = 8023.799999999999
$ws.Range("M136").Value = -107143680
$ws.Range("N136").Value = -13123.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 30231918
$ws.Range("I121").Value = 1052.5
$ws.Range("K121").Value = 3157.5
$ws.Range("M121").Value = -1847.5
$ws.Range("H132").Value = 1080.0834
$ws.Range("I132").Value = 730.7059
$ws.Range("K132").Value = 6576.3531
$ws.Range("M132").Value = -4046.3531

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2281.8667
$ws.Range("I126").Value = 1966.6666
$ws.Range("J126").Value = 2754.6667
$ws.Range("K126").Value = 5899.9998
$ws.Range("L126").Value = 8264.000100000001
$ws.Range("M126").Value = -3429.9998
$ws.Range("N126").Value = -13204.0001
$ws.Range("H132").Value = 75899.07000000001
$ws.Range("I132").Value = 91807.63
$ws.Range("J132").Value = 64961.938
$ws.Range("K132").Value = 275422.89
$ws.Range("L132").Value = 194885.814
$ws.Range("M132").Value = -272892.89
$ws.Range("N132").Value = -199945.814

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4227.636
$ws.Range("I40").Value = 2959.8
$ws.Range("J40").Value = 5284.1665
$ws.Range("K40").Value = 2959.8
$ws.Range("L40").Value = 5284.1665
$ws.Range("M40").Value = -2823.8
$ws.Range("N40").Value = -5556.1665
$ws.Range("H132").Value = 34771.387
$ws.Range("I132").Value = 1345.875
$ws.Range("J132").Value = 70425.266
$ws.Range("K132").Value = 4037.625
$ws.Range("L132").Value = 211275.798
$ws.Range("M132").Value = -1507.625
$ws.Range("N132").Value = -216335.798
$ws.Range("H136").Value = 126248.75
$ws.Range("I136").Value = 143911.42
$ws.Range("J136").Value = 112511.11
$ws.Range("K136").Value = 431734.26
$ws.Range("L136").Value = 337533.33
$ws.Range("M136").Value = -429184.26
$ws.Range("N136").Value = -337533.33

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 63745.562
$ws.Range("I132").Value = 39400
$ws.Range("J132").Value = 169243
$ws.Range("K132").Value = 118200
$ws.Range("L132").Value = 507729
$ws.Range("M132").Value = -115670
$ws.Range("N132").Value = -512789
$ws.Range("H136").Value = 66801.67999999999
$ws.Range("I136").Value = 43018.668
$ws.Range("K136").Value = 129056.004
$ws.Range("M136").Value = -126506.004
